$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: run a scoped Find & Replace (no wildcards, case-sensitive) limited
# to a single paragraph's Range so it can never bleed into a neighboring,
# textually-similar paragraph elsewhere in the resume.
# ---------------------------------------------------------------------------
function Replace-InParagraph($para, [string]$find, [string]$replace) {
    $rng = $para.Range
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 0, $false, $replace, 2)
    return $ok
}

# ---------------------------------------------------------------------------
# 1) Podhub bullet:
#    "Node.js, Express.js, ReactJS, Bootstrap framework, API calls/routing."
#    -> "Express.js, ReactJS, Node.js, Bootstrap, Listen Notes API."
#    Split the edit into two scoped replaces that each avoid touching the
#    "ReactJS" run so its spell-check (proofErr) wrapping is left intact.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ReactJS, Bootstrap framework, API calls/routing*") {
        Replace-InParagraph $p "Node.js, Express.js, " "Express.js, " | Out-Null
        Replace-InParagraph $p ", Bootstrap framework, API calls/routing." ", Node.js, Bootstrap, Listen Notes API." | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Mammoth bullet:
#    "Node.js, Express.js, Google Maps APIs, MySQL, Handlebars, CSS."
#    -> "Express.js, Node.js, MySQL, Google Maps API."
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Google Maps APIs*") {
        Replace-InParagraph $p "Node.js, Express.js, Google Maps APIs, MySQL, Handlebars, CSS." "Express.js, Node.js, MySQL, Google Maps API." | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Recipe Finder bullet:
#    "Node.js, Express.js, MySQL, Bootstrap, Handlebars, API calls/routing. "
#    -> "Express.js, Node.js, MySQL, Bootstrap, Yummly API. "
#    Again split around "Bootstr"+"ap" so that those two untouched runs (and
#    the new "Yummly" spell-check run) stay separate from the edited text.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*MySQL, Bootstrap, Handlebars*") {
        Replace-InParagraph $p "Node.js, Express.js, MySQL, Bootstr" "Express.js, Node.js, MySQL, Bootstr" | Out-Null
        Replace-InParagraph $p ", Handlebars, API calls/routing." ", Yummly API." | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Move the hidden "_GoBack" bookmark from the end of the last bullet in
#    the Work Experience section to the end of the Sudoku project's URL
#    line. A collapsed Range can't be handed to Bookmarks.Add directly in
#    this host, so we temporarily insert a one-character placeholder, wrap
#    the bookmark around it, then delete the placeholder again - leaving the
#    bookmark collapsed exactly where it should be.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*sudoku-cjy.herokuapp.com*") {
        $endPos = $p.Range.End - 1
        $insertion = $d.Range($endPos, $endPos)
        $insertion.InsertAfter("X")
        $marker = $d.Range($endPos, $endPos + 1)
        $d.Bookmarks.Add("_GoBack", $marker)
        $d.Range($endPos, $endPos + 1).Text = ""
        break
    }
}

Write-Output "done"
